$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999616220903,
    0.99385465086959945,
    0.9705951338219827,
    0.95932253566473369,
    0.94841558446867635,
    0.92192430631314026,
    0.91532038747726419,
    0.90618172526226748,
    0.89442486940946897,
    0.88358862122698134,
    0.88199157713947485,
    0.87921025209988135,
    0.86791808166936291,
    0.86374957057094615,
    0.8611574945962357,
    0.85865059278242906,
    0.85494224446020795,
    0.85383329852041312,
    0.99207550873720396,
    0.97847723725522229,
    0.97707871981300887,
    0.97581420879155256,
    0.9635592952092259,
    0.95053736565312308,
    0.94408018925236981,
    0.92003212954000435,
    0.91764487755063151,
    0.90723332250367172,
    0.90019696849852493,
    0.89787105634329145,
    0.90062798365726815,
    0.89894864066997182,
    0.89842862309220528
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
